# 57 - ERD maken
# Update the "Users" entity attribute list on the ClassDiagram/EntityDraft
# sheet: several entity id columns are renamed to the generic "id", the
# Users entity gains id/firstname/lastname/department/email fields (and
# keeps password/salt/userid below them), and the Subscribes "DepartmentId"
# field becomes "Department".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Entity id headers (row 5) that become the generic "id" field name.
$ws.Range("A5").Value = "id"
$ws.Range("B5").Value = "id"
$ws.Range("D5").Value = "id"
$ws.Range("E5").Value = "id"
$ws.Range("G5").Value = "id"
$ws.Range("H5").Value = "id"

# Users entity (column A) gains new fields; existing ones shift down.
$ws.Range("A6").Value = "firstname"
$ws.Range("A7").Value = "lastname"
$ws.Range("A8").Value = "department"
$ws.Range("A9").Value = "email"
$ws.Range("A10").Value = "password"
$ws.Range("A11").Value = "salt"
$ws.Range("A12").Value = "userid"

# WantedSubscribers entity: DepartmentId -> Department
$ws.Range("J7").Value = "Department"

# Update the current selection to A9 (and reset the scrolled top-left cell).
$ws.Range("A9").Select()
